# Apply updated cryptocurrency price/volume data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "Price" (column D) cell to an exact text string, preserving
# the General/default cell style (avoids Excel auto-converting values like
# "57.394.00" or "0.0900" into numbers and losing formatting / trailing zeros).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "57.394.00"
$ws.Range("E2").Value = "  -2.61%  "
Set-TextValue $ws.Range("D3") "2.419.56"
$ws.Range("E3").Value = "  -3.55%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "512.36"
$ws.Range("E5").Value = "  -3.75%  "
Set-TextValue $ws.Range("D6") "128.73"
$ws.Range("E6").Value = "  -4.76%  "
Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  -0.13%  "
Set-TextValue $ws.Range("D8") "0.548"
$ws.Range("E8").Value = "  -3.16%  "
Set-TextValue $ws.Range("D9") "2.426.73"
$ws.Range("E9").Value = "  -3.41%  "
Set-TextValue $ws.Range("D10") "0.0956"
$ws.Range("E10").Value = "  -5.67%  "
$ws.Range("E11").Value = "  -1.58%  "
Set-TextValue $ws.Range("D12") "5.19"
$ws.Range("E12").Value = "  -3.97%  "
Set-TextValue $ws.Range("D13") "0.331"
$ws.Range("E13").Value = "  -3.90%  "
Set-TextValue $ws.Range("D14") "2.847.26"
$ws.Range("E14").Value = "  -3.64%  "
Set-TextValue $ws.Range("D15") "57.323.97"
$ws.Range("E15").Value = "  -2.57%  "
Set-TextValue $ws.Range("D16") "21.41"
$ws.Range("E16").Value = "  -5.80%  "
$ws.Range("E17").Value = "  -4.41%  "
Set-TextValue $ws.Range("D18") "2.416.76"
$ws.Range("E18").Value = "  -3.77%  "
Set-TextValue $ws.Range("D19") "10.35"
$ws.Range("E19").Value = "  -5.97%  "
Set-TextValue $ws.Range("D20") "314.03"
$ws.Range("E20").Value = "  -2.56%  "
Set-TextValue $ws.Range("D21") "4.08"
$ws.Range("E21").Value = "  -4.14%  "
$ws.Range("E22").Value = "  -0.09%  "
Set-TextValue $ws.Range("D23") "5.64"
$ws.Range("E23").Value = "  -4.97%  "
Set-TextValue $ws.Range("D24") "63.51"
$ws.Range("E24").Value = "  -2.27%  "
Set-TextValue $ws.Range("D25") "0.401"
$ws.Range("E25").Value = "  -4.30%  "
Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -2.59%  "
Set-TextValue $ws.Range("D28") "7.19"
$ws.Range("E28").Value = "  -4.80%  "
Set-TextValue $ws.Range("D29") "169.25"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -6.00%  "
$ws.Range("E31").Value = "  -5.17%  "
Set-TextValue $ws.Range("D32") "6.16"
$ws.Range("E32").Value = "  -4.90%  "
Set-TextValue $ws.Range("D33") "1.15"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -4.09%  "
Set-TextValue $ws.Range("D37") "1.27"
$ws.Range("E37").Value = "  -7.59%  "
$ws.Range("E38").Value = "  -4.63%  "
Set-TextValue $ws.Range("D39") "36.32"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("E40").Value = "  -5.53%  "
Set-TextValue $ws.Range("D41") "0.766"
$ws.Range("E41").Value = "  -4.15%  "
Set-TextValue $ws.Range("D42") "3.36"
$ws.Range("E42").Value = "  -5.87%  "
Set-TextValue $ws.Range("D43") "264.74"
$ws.Range("E43").Value = "  -5.80%  "
Set-TextValue $ws.Range("D44") "4.87"
$ws.Range("E44").Value = "  -2.59%  "
Set-TextValue $ws.Range("D45") "0.581"
$ws.Range("E45").Value = "  -3.72%  "
Set-TextValue $ws.Range("D46") "122.26"
$ws.Range("E46").Value = "  -5.71%  "
Set-TextValue $ws.Range("D47") "0.0900"
$ws.Range("E47").Value = "  -2.74%  "
Set-TextValue $ws.Range("D48") "0.0480"
$ws.Range("E48").Value = "  -3.82%  "
Set-TextValue $ws.Range("D49") "0.0210"
$ws.Range("E49").Value = "  -3.67%  "
Set-TextValue $ws.Range("D50") "16.46"
$ws.Range("E50").Value = "  -4.48%  "
Set-TextValue $ws.Range("D51") "1.689.20"
$ws.Range("E51").Value = "  -3.93%  "
